$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row right-answer count
$ws.Range("B11").Value = 5

# Update "Total" row right-answer count and the Corr/total marks summary
$ws.Range("B12").Value = 125
$ws.Range("E12").Value = "125/140"
